$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.761.18"
$ws.Range("E2").Value = "  -4.68%  "
$ws.Range("D3").Value = "2.431.53"
$ws.Range("E3").Value = "  -4.29%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'308.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").Value = "'93.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.75%  "
$ws.Range("D7").Value = "'0.547"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.31%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -5.85%  "
$ws.Range("D10").Value = "'33.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.88%  "
$ws.Range("E11").Value = "  -3.66%  "
$ws.Range("E12").Value = "  -0.78%  "
$ws.Range("E13").Value = "  -6.40%  "
$ws.Range("D14").Value = "2.811.43"
$ws.Range("E14").Value = "  -3.82%  "
$ws.Range("D15").Value = "2.434.79"
$ws.Range("E15").Value = "  -5.24%  "
$ws.Range("E16").Value = "  -9.86%  "
$ws.Range("D17").Value = "'0.779"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.27%  "
$ws.Range("D18").Value = "40.798.78"
$ws.Range("E18").Value = "  -4.54%  "
$ws.Range("E19").Value = "  -7.34%  "
$ws.Range("E20").Value = "  -5.14%  "
$ws.Range("D21").Value = "'11.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.31%  "
$ws.Range("D22").Value = "'66.78"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.74%  "
$ws.Range("D23").Value = "'235.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.75%  "
$ws.Range("E24").Value = "  -5.27%  "
$ws.Range("E25").Value = "  -6.77%  "
$ws.Range("E26").Value = "  +6.84%  "
$ws.Range("D27").Value = "'24.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.89%  "
$ws.Range("D28").Value = "'2.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.48%  "
$ws.Range("E29").Value = "  -5.79%  "
$ws.Range("D30").Value = "'35.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -9.09%  "
$ws.Range("D31").Value = "'152.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.44%  "
$ws.Range("D32").Value = "'5.51"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.89%  "
$ws.Range("E33").Value = "  -0.83%  "
$ws.Range("D34").Value = "'2.52"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.81%  "
$ws.Range("E35").Value = "  -6.30%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "'2.98"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.40%  "
$ws.Range("B37").Value = "Celestia"
$ws.Range("C37").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D37").Value = "'17.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.19%  "
$ws.Range("E38").Value = "  -8.79%  "
$ws.Range("E39").Value = "  -5.18%  "
$ws.Range("E40").Value = "  -9.35%  "
$ws.Range("D41").Value = "'4.07"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.89%  "
$ws.Range("D42").Value = "'21.12"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.08%  "
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("D44").Value = "1.955.00"
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("E45").Value = "  -6.06%  "
$ws.Range("E46").Value = "  -9.06%  "
$ws.Range("E47").Value = "  -3.05%  "
$ws.Range("D48").Value = "'76.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.88%  "
$ws.Range("D49").Value = "'96.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'68.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.32%  "
$ws.Range("D51").Value = "'0.177"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.98%  "
